$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.52698182277397
$ws.Range("C2").Value = 3.023968074486815
$ws.Range("D2").Value = 3.904400907923373
$ws.Range("E2").Value = 10.77759118223912
$ws.Range("F2").Value = 69.53616316249607
$ws.Range("J2").Value = 11.27624633093274
$ws.Range("K2").Value = 17.684123285992
$ws.Range("L2").Value = 11.10012192625151
$ws.Range("B3").Value = 21.54617254014067
$ws.Range("C3").Value = 2.974724868052011
$ws.Range("D3").Value = 3.78020419545809
$ws.Range("E3").Value = 10.81090140214611
$ws.Range("F3").Value = 68.71104351508627
$ws.Range("J3").Value = 11.26839204990302
$ws.Range("K3").Value = 17.71754574697328
$ws.Range("L3").Value = 11.14497077895334
$ws.Range("B4").Value = 21.56471972667945
$ws.Range("C4").Value = 2.946699213255007
$ws.Range("D4").Value = 3.700840707651334
$ws.Range("E4").Value = 10.83268889960536
$ws.Range("F4").Value = 68.20216886422497
$ws.Range("J4").Value = 11.26399002304225
$ws.Range("K4").Value = 17.74346305139033
$ws.Range("L4").Value = 11.17448627261466
$ws.Range("B5").Value = 21.57397614888351
$ws.Range("C5").Value = 2.935857803289229
$ws.Range("D5").Value = 3.667736998744724
$ws.Range("E5").Value = 10.84190398344518
$ws.Range("F5").Value = 67.9943717979892
$ws.Range("J5").Value = 11.26230208162493
$ws.Range("K5").Value = 17.75537804697906
$ws.Range("L5").Value = 11.18701203763181
$ws.Range("B6").Value = 21.57561564859268
$ws.Range("C6").Value = 2.934093227693755
$ws.Range("D6").Value = 3.662194687350417
$ws.Range("E6").Value = 10.84345448954108
$ws.Range("F6").Value = 67.95984553633522
$ws.Range("J6").Value = 11.26202820225961
$ws.Range("K6").Value = 17.7574381618063
$ws.Range("L6").Value = 11.18912202337141
$ws.Range("B7").Value = 21.56483768992824
$ws.Range("C7").Value = 2.946550627245458
$ws.Range("D7").Value = 3.700397320413239
$ws.Range("E7").Value = 10.83281181386298
$ws.Range("F7").Value = 68.19936798461745
$ws.Range("J7").Value = 11.26396682978021
$ws.Range("K7").Value = 17.74361826586593
$ws.Range("L7").Value = 11.17465318251982
$ws.Range("B8").Value = 21.53219366327421
$ws.Range("C8").Value = 3.006545158122091
$ws.Range("D8").Value = 3.862233136271215
$ws.Range("E8").Value = 10.78879999073006
$ws.Range("F8").Value = 69.25221172474738
$ws.Range("J8").Value = 11.2734509998202
$ws.Range("K8").Value = 17.69452561918763
$ws.Range("L8").Value = 11.11517557995327
$ws.Range("B9").Value = 21.52191954814831
$ws.Range("C9").Value = 3.203828611161989
$ws.Range("D9").Value = 4.154229233018819
$ws.Range("E9").Value = 10.71304731982896
$ws.Range("F9").Value = 71.29292669333675
$ws.Range("J9").Value = 11.29537928146185
$ws.Range("K9").Value = 17.64119333070734
$ws.Range("L9").Value = 11.01421147073383
$ws.Range("B10").Value = 21.54716606557569
$ws.Range("C10").Value = 3.445411907779055
$ws.Range("D10").Value = 4.352469972263178
$ws.Range("E10").Value = 10.66377474212581
$ws.Range("F10").Value = 72.76961114860076
$ws.Range("J10").Value = 11.31351182643871
$ws.Range("K10").Value = 17.62832200886974
$ws.Range("L10").Value = 10.9495541501448
$ws.Range("B11").Value = 21.56575678325895
$ws.Range("C11").Value = 3.549603523146473
$ws.Range("D11").Value = 4.43897720641746
$ws.Range("E11").Value = 10.64273455817701
$ws.Range("F11").Value = 73.43480147946708
$ws.Range("J11").Value = 11.32219778337064
$ws.Range("K11").Value = 17.62819591828723
$ws.Range("L11").Value = 10.92220065632664
$ws.Range("B12").Value = 21.57381518127467
$ws.Range("C12").Value = 3.58823716231884
$ws.Range("D12").Value = 4.471197135727372
$ws.Range("E12").Value = 10.63496396522376
$ws.Range("F12").Value = 73.68561018233595
$ws.Range("J12").Value = 11.32554956237772
$ws.Range("K12").Value = 17.62897218307058
$ws.Range("L12").Value = 10.91213828440944
$ws.Range("B13").Value = 21.5720344284445
$ws.Range("C13").Value = 3.579953257184008
$ws.Range("D13").Value = 4.464282095176129
$ws.Range("E13").Value = 10.63662875728743
$ws.Range("F13").Value = 73.63164448069602
$ws.Range("J13").Value = 11.32482491782296
$ws.Range("K13").Value = 17.62876835703995
$ws.Range("L13").Value = 10.91429224562151
$ws.Range("B14").Value = 21.56639936350165
$ws.Range("C14").Value = 3.552798401456842
$ws.Range("D14").Value = 4.441638805175279
$ws.Range("E14").Value = 10.64209132546415
$ws.Range("F14").Value = 73.4554581919451
$ws.Range("J14").Value = 11.32247228352146
$ws.Range("K14").Value = 17.62824326911403
$ws.Range("L14").Value = 10.92136689251632
$ws.Range("B15").Value = 21.56308024331064
$ws.Range("C15").Value = 3.536058276795824
$ws.Range("D15").Value = 4.427698729416954
$ws.Range("E15").Value = 10.64546292259092
$ws.Range("F15").Value = 73.34739364873602
$ws.Range("J15").Value = 11.32103936966233
$ws.Range("K15").Value = 17.62802894155919
$ws.Range("L15").Value = 10.92573883033872
$ws.Range("B16").Value = 21.54609385001104
$ws.Range("C16").Value = 3.438487819279214
$ws.Range("D16").Value = 4.346741694147598
$ws.Range("E16").Value = 10.6651773562474
$ws.Range("F16").Value = 72.72599524839245
$ws.Range("J16").Value = 11.31295295480439
$ws.Range("K16").Value = 17.62844556435524
$ws.Range("L16").Value = 10.95138317687612
$ws.Range("B17").Value = 21.53749124913212
$ws.Range("C17").Value = 3.377169094787807
$ws.Range("D17").Value = 4.296128578311351
$ws.Range("E17").Value = 10.67762295251361
$ws.Range("F17").Value = 72.34300574043714
$ws.Range("J17").Value = 11.30810401767062
$ws.Range("K17").Value = 17.63016885814894
$ws.Range("L17").Value = 10.96764237662625
$ws.Range("B18").Value = 21.53321237915243
$ws.Range("C18").Value = 3.34136311495288
$ws.Range("D18").Value = 4.266672007709166
$ws.Range("E18").Value = 10.68491071602453
$ws.Range("F18").Value = 72.12211101564787
$ws.Range("J18").Value = 11.30535622398132
$ws.Range("K18").Value = 17.63169933834636
$ws.Range("L18").Value = 10.9771881025609
$ws.Range("B19").Value = 21.5318786315275
$ws.Range("C19").Value = 3.329147628686835
$ws.Range("D19").Value = 4.256639557868133
$ws.Range("E19").Value = 10.68740047215121
$ws.Range("F19").Value = 72.04721955997037
$ws.Range("J19").Value = 11.30443295081852
$ws.Range("K19").Value = 17.63231013919669
$ws.Range("L19").Value = 10.98045342843016
$ws.Range("B20").Value = 21.53833778010338
$ws.Range("C20").Value = 3.383752172526731
$ws.Range("D20").Value = 4.301552240909615
$ws.Range("E20").Value = 10.67628471199155
$ws.Range("F20").Value = 72.38383971842883
$ws.Range("J20").Value = 11.30861593470956
$ws.Range("K20").Value = 17.62992959314755
$ws.Range("L20").Value = 10.96589149544308
$ws.Range("B21").Value = 21.56802690921882
$ws.Range("C21").Value = 3.560796738178586
$ws.Range("D21").Value = 4.448304377423229
$ws.Range("E21").Value = 10.64048149964096
$ws.Range("F21").Value = 73.50723891427002
$ws.Range("J21").Value = 11.32316161290166
$ws.Range("K21").Value = 17.62837513878256
$ws.Range("L21").Value = 10.91928087188507
$ws.Range("B22").Value = 21.59336469153131
$ws.Range("C22").Value = 3.671717658727021
$ws.Range("D22").Value = 4.5410728859113
$ws.Range("E22").Value = 10.61822919393141
$ws.Range("F22").Value = 74.2350613434123
$ws.Range("J22").Value = 11.33303275587648
$ws.Range("K22").Value = 17.63216177656185
$ws.Range("L22").Value = 10.89054215134579
$ws.Range("B23").Value = 21.57929985661788
$ws.Range("C23").Value = 3.612955255812246
$ws.Range("D23").Value = 4.49185115990368
$ws.Range("E23").Value = 10.63000094307968
$ws.Range("F23").Value = 73.8472386378955
$ws.Range("J23").Value = 11.32773107501592
$ws.Range("K23").Value = 17.62970146513658
$ws.Range("L23").Value = 10.90572290868888
$ws.Range("B24").Value = 21.53795298593662
$ws.Range("C24").Value = 3.380777685609757
$ws.Range("D24").Value = 4.299101318705677
$ws.Range("E24").Value = 10.67688931736617
$ws.Range("F24").Value = 72.36538087437962
$ws.Range("J24").Value = 11.30838437264594
$ws.Range("K24").Value = 17.63003608363799
$ws.Range("L24").Value = 10.96668245174564
$ws.Range("B25").Value = 21.5189378834121
$ws.Range("C25").Value = 3.109816227801651
$ws.Range("D25").Value = 4.078043950034949
$ws.Range("E25").Value = 10.7324158429083
$ws.Range("F25").Value = 70.74425514754192
$ws.Range("J25").Value = 11.28909136799111
$ws.Range("K25").Value = 17.65100780804013
$ws.Range("L25").Value = 11.03985065355978
